# Auto-generated update script for Zalera_Profits workbook
# Updates currentAveragePrice / Leve price / profit columns (H-N) for select rows
# across multiple job-class sheets, reflecting refreshed market-board data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$values_ALC = @{
    "H62" = 66668844
    "I62" = 66668844
    "J62" = 0
    "K62" = 66668844
    "L62" = 0
    "M62" = -66668220
    "H64" = 4881
    "I64" = 5008.3335
    "J64" = 4499
    "K64" = 5008.3335
    "L64" = 4499
    "M64" = -4760.3335
    "N64" = -4995
    "H65" = 66668844
    "I65" = 66668844
    "J65" = 0
    "K65" = 333344220
    "L65" = 0
    "M65" = -333341100
    "H67" = 4881
    "I67" = 5008.3335
    "J67" = 4499
    "K67" = 5008.3335
    "L67" = 4499
    "M67" = -4150.3335
    "N67" = -6215
    "H76" = 6671354.5
    "I76" = 11115258
    "J76" = 5499.1665
    "K76" = 11115258
    "L76" = 5499.1665
    "M76" = -11114943
    "N76" = -6129.1665
    "H79" = 6671354.5
    "I79" = 11115258
    "J79" = 5499.1665
    "K79" = 11115258
    "L79" = 5499.1665
    "M79" = -11114166
    "N79" = -7683.1665
    "H132" = 2430.6086
    "I132" = 1809.5454
    "K132" = 5428.6362
    "M132" = -2898.6362
    "H135" = 4173409.5
    "J135" = 19983.715
    "L135" = 179853.435
    "N135" = -184923.435
    "H137" = 7208.6978
    "I137" = 3987.0588
    "J137" = 9315.154
    "K137" = 11961.1764
    "L137" = 27945.462
    "M137" = -9411.1764
    "N137" = -33045.462
    "H138" = 4645.973
    "I138" = 4130.8184
    "J138" = 4863.923
    "K138" = 12392.4552
    "L138" = 14591.769
    "M138" = -7252.4552
    "N138" = -24871.769
    "H141" = 1934.1666
    "I141" = 962.2222
    "J141" = 4850
    "K141" = 2886.6666
    "L141" = 14550
    "M141" = 2293.3334
    "N141" = -24910
}
foreach ($addr in $values_ALC.Keys) {
    $ws.Range($addr).Value = $values_ALC[$addr]
}
$clears_ALC = @("N62","N65")
foreach ($addr in $clears_ALC) {
    $ws.Range($addr).ClearContents()
}

$ws = $wb.Worksheets.Item("ARM")
$values_ARM = @{
    "H2" = 8367317.5
    "I2" = 8367317.5
    "J2" = 0
    "K2" = 8367317.5
    "L2" = 0
    "M2" = -8367204.5
    "H45" = 3132.4
    "J45" = 4250
    "L45" = 4250
    "N45" = -5004
    "H116" = 8367317.5
    "I116" = 8367317.5
    "J116" = 0
    "K116" = 8367317.5
    "L116" = 0
    "M116" = -8365023.5
}
foreach ($addr in $values_ARM.Keys) {
    $ws.Range($addr).Value = $values_ARM[$addr]
}
$clears_ARM = @("N2","N116")
foreach ($addr in $clears_ARM) {
    $ws.Range($addr).ClearContents()
}

$ws = $wb.Worksheets.Item("BSM")
$values_BSM = @{
    "H3" = 8367317.5
    "I3" = 8367317.5
    "J3" = 0
    "K3" = 8367317.5
    "L3" = 0
    "M3" = -8367203.5
    "H134" = 5034.5835
    "I134" = 3277.84
    "J134" = 9027.182000000001
    "K134" = 9833.52
    "L134" = 27081.546
    "M134" = -7298.52
    "N134" = -32151.546
}
foreach ($addr in $values_BSM.Keys) {
    $ws.Range($addr).Value = $values_BSM[$addr]
}
$clears_BSM = @("N3")
foreach ($addr in $clears_BSM) {
    $ws.Range($addr).ClearContents()
}

$ws = $wb.Worksheets.Item("CRP")
$values_CRP = @{
    "H31" = 2990.2708
    "I31" = 1311.9642
    "J31" = 5339.9
    "K31" = 1311.9642
    "L31" = 5339.9
    "M31" = -1016.9642
    "N31" = -5929.9
    "H34" = 2990.2708
    "I34" = 1311.9642
    "J34" = 5339.9
    "K34" = 1311.9642
    "L34" = 5339.9
    "M34" = -1109.9642
    "N34" = -5743.9
    "H132" = 99905.08
    "I132" = 5596.375
    "J132" = 250799
    "K132" = 16789.125
    "L132" = 752397
    "M132" = -14259.125
    "N132" = -757457
}
foreach ($addr in $values_CRP.Keys) {
    $ws.Range($addr).Value = $values_CRP[$addr]
}

$ws = $wb.Worksheets.Item("CUL")
$values_CUL = @{
    "H4" = 6603407.5
    "I4" = 9278283
    "K4" = 27834849
    "M4" = -27834737
    "H44" = 3377.2222
    "I44" = 950
    "J44" = 4070.7144
    "K44" = 2850
    "L44" = 12212.1432
    "M44" = -2452
    "N44" = -13008.1432
    "H68" = 83469.58
    "I68" = 288159.72
    "J68" = 8057.421
    "K68" = 864479.1599999999
    "L68" = 24172.263
    "M68" = -863668.1599999999
    "N68" = -25794.263
    "H71" = 83469.58
    "I71" = 288159.72
    "J71" = 8057.421
    "K71" = 2593437.48
    "L71" = 72516.789
    "M71" = -2589381.48
    "N71" = -80628.789
}
foreach ($addr in $values_CUL.Keys) {
    $ws.Range($addr).Value = $values_CUL[$addr]
}

$ws = $wb.Worksheets.Item("GSM")
$values_GSM = @{
    "H39" = 31704.2
    "I39" = 27261
    "J39" = 32815
    "K39" = 27261
    "L39" = 32815
    "N39" = -33879
    "M39" = -26729
    "H51" = 96117.64999999999
    "J51" = 96117.64999999999
    "L51" = 96117.64999999999
    "N51" = -97135.64999999999
    "H70" = 9340.559999999999
    "I70" = 9751.723
    "J70" = 8283.286
    "K70" = 9751.723
    "L70" = 8283.286
    "M70" = -9481.723
    "N70" = -8823.286
    "H73" = 9340.559999999999
    "I73" = 9751.723
    "J73" = 8283.286
    "K73" = 9751.723
    "L73" = 8283.286
    "M73" = -8815.723
    "N73" = -10155.286
    "H97" = 677.8421
    "I97" = 692.26666
    "J97" = 623.75
    "K97" = 692.26666
    "L97" = 623.75
    "M97" = -196.26666
    "N97" = -1615.75
    "H102" = 1563.6666
    "I102" = 1609.7742
    "J102" = 1277.8
    "K102" = 1609.7742
    "L102" = 1277.8
    "M102" = 12.22579999999994
    "N102" = -4521.8
    "H126" = 3170.15
    "I126" = 2235.6155
    "J126" = 4905.7144
    "K126" = 6706.8465
    "L126" = 14717.1432
    "M126" = -4236.8465
    "N126" = -19657.1432
}
foreach ($addr in $values_GSM.Keys) {
    $ws.Range($addr).Value = $values_GSM[$addr]
}

$ws = $wb.Worksheets.Item("LTW")
$values_LTW = @{
    "H68" = 1948.75
    "I68" = 2098.3333
    "K68" = 2098.3333
    "M68" = -1349.3333
    "H71" = 1948.75
    "I71" = 2098.3333
    "K71" = 10491.6665
    "M71" = -6747.666499999999
    "H82" = 1491.9
    "I82" = 1371.8182
    "J82" = 1561.421
    "K82" = 1371.8182
    "L82" = 1561.421
    "M82" = -1010.8182
    "N82" = -2283.421
    "H85" = 1491.9
    "I85" = 1371.8182
    "J85" = 1561.421
    "K85" = 1371.8182
    "L85" = 1561.421
    "M85" = -123.8181999999999
    "N85" = -4057.421
}
foreach ($addr in $values_LTW.Keys) {
    $ws.Range($addr).Value = $values_LTW[$addr]
}

$ws = $wb.Worksheets.Item("WVR")
$values_WVR = @{
    "H63" = 18734.637
    "I63" = 20833.166
    "J63" = 16216.4
    "K63" = 20833.166
    "L63" = 16216.4
    "M63" = -20209.166
    "N63" = -17464.4
    "H66" = 18734.637
    "I66" = 20833.166
    "J66" = 16216.4
    "K66" = 62499.49800000001
    "L66" = 48649.2
    "M66" = -59379.49800000001
    "N66" = -54889.2
    "H132" = 5907.4253
    "I132" = 4555
    "J132" = 8294.058999999999
    "K132" = 13665
    "L132" = 24882.177
    "M132" = -11135
    "N132" = -29942.177
}
foreach ($addr in $values_WVR.Keys) {
    $ws.Range($addr).Value = $values_WVR[$addr]
}
